# Sync attendance_reports: reorder "Recorded By" (column G) entries.
# Applies a fixed set of text replacements to the "Recorded By" column,
# swapping the order of the first two comma-separated recorder names
# while leaving the rest of the list (and any unrelated values) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "dnasr281@gmail.com, System"                = "System, dnasr281@gmail.com"
    "System, admin@admin.com"                   = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com"       = "admin@admin.com, dnasr281@gmail.com"
    "system, backup@backdoor.com, System"       = "backup@backdoor.com, system, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
